$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '67.212.23'
$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  +0.45%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.624.54'
$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  +0.74%  '
$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '597.08'
$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  +0.97%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '152.38'
$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  -0.72%  '
$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.554'
$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  +2.91%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '2.623.47'
$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  +0.80%  '
$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  +0.73%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '5.18'
$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  -0.31%  '
$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  -1.05%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '27.54'
$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  +0.83%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '3.102.50'
$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  +0.84%  '
$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  +0.62%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '67.168.01'
$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  +0.43%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '2.627.45'
$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  +0.97%  '
$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  -0.27%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '363.34'
$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  +2.32%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '7.48'
$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  -3.04%  '
$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  -0.31%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '2.10'
$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +3.52%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  +0.05%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '70.93'
$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  +7.05%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '10.01'
$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  -2.09%  '
$c = $ws.Range('B27')
$c.NumberFormat = "@"
$c.Value = 'WrappedeETH'
$c = $ws.Range('C27')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '2.760.08'
$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  +0.61%  '
$c = $ws.Range('B28')
$c.NumberFormat = "@"
$c.Value = 'Binance-PegBSC-USD'
$c = $ws.Range('C28')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '1.02'
$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  +1.80%  '
$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  +0.50%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '574.75'
$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  -5.61%  '
$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -3.30%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '7.81'
$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -1.53%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '1.83'
$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  -0.62%  '
$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  +0.09%  '
$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  -3.61%  '
$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  -1.46%  '
$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -1.08%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '157.54'
$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  +1.75%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '19.17'
$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  +0.25%  '
$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  -2.55%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '1.81'
$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  +0.21%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.55'
$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  +0.88%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '41.19'
$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  -0.19%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '16.35'
$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  -0.55%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '156.10'
$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  +0.97%  '
$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  -2.03%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '3.73'
$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '
$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  -0.22%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '20.57'
$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  -0.56%  '
